# Summer 2019 word-count entry ("Happenin's 'Round the House", row 9) was
# forgotten/mis-entered; correct the actual-delivered figures and bring the
# selection to rest on the next row, matching how the author would have
# left the sheet after fixing the numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 9: "Actual Delivered" words (F9) and pages (G9)
$ws.Range("F9").Value = 1421
$ws.Range("G9").Value = 7

# Make sure the Totals/Pages-per-section formulas (F11:H12) are
# recalculated from the corrected figures.
$excel.Calculate()

# Leave the selection where the author ended up after the edit.
$ws.Range("G10").Select()
